# Automatische test-sync: 2025-07-27 18:37:50
# Adds a new test-mail log entry (row 5) to the "Logs" sheet and
# updates the related conditional formatting ranges + dashboard count.

$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append the new row with the 4th test mail ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(5, 1).Value = "Kun je dit intern bespreken?"
$logs.Cells.Item(5, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(5, 3).Value = "Testmail #4: Kun je dit intern bespreken?"
$logs.Cells.Item(5, 4).Value = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item(5, 5).Value = "Beste afzender,`r`nBedankt voor je e-mail. Kun je wat meer specifieke informatie geven over waarover je precies wilt dat er intern overlegd wordt? Op die manier kan ik ervoor zorgen dat je aanvraag bij de juiste persoon of afdeling terechtkomt.`r`nMet vriendelijke groet,`r`n[Naam]`r`nE-mailassistent"
$logs.Cells.Item(5, 6).Value = "2025-07-27 18:37:05"
$logs.Cells.Item(5, 7).Value = "Ja"
$logs.Cells.Item(5, 8).Value = "Nee"
$logs.Cells.Item(5, 9).Value = "Ja"
$logs.Cells.Item(5, 10).Value = "Ja"

# --- Extend the conditional formatting ranges from row 4 to row 5 ---
$logs.Range("D2:D4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D5"))
$logs.Range("G2:G4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G5"))
$logs.Range("H2:H4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H5"))
$logs.Range("I2:I4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I5"))
$logs.Range("J2:J4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J5"))

# --- Sheet "Dashboard": bump the "Intern verzoek / Actie voor medewerker" count ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(3, 2).Value = 2
